# Edit script: applies the re-ordering of betting-odds rows and appends
# two new match rows (91-92), matching the upstream scraper re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-synchronised match rows: the scraper re-ordered several fixtures
#     played on the same matchday, so columns F:V (everything except the
#     fixed Indice/pais/torneio/temporada/data_partida columns A:E) need to
#     be rewritten with the data belonging to the new row order. ---

# Row 6 now holds the fixture previously shown on row 8
$ws.Cells.Item(6, 6).Value = "Nantes"
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = "Toulouse"
$ws.Cells.Item(6, 9).Value = 2
$ws.Cells.Item(6, 10).Value = 2.27
$ws.Cells.Item(6, 11).Value = "05/07/2023 22:22"
$ws.Cells.Item(6, 12).Value = 3.09
$ws.Cells.Item(6, 13).Value = "13/08/2023 14:26"
$ws.Cells.Item(6, 14).Value = 3.38
$ws.Cells.Item(6, 15).Value = "05/07/2023 22:22"
$ws.Cells.Item(6, 16).Value = 3.48
$ws.Cells.Item(6, 17).Value = "13/08/2023 14:26"
$ws.Cells.Item(6, 18).Value = 3.37
$ws.Cells.Item(6, 19).Value = "05/07/2023 22:22"
$ws.Cells.Item(6, 20).Value = 2.43
$ws.Cells.Item(6, 21).Value = "13/08/2023 14:26"
$ws.Cells.Item(6, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/nantes-toulouse/A1jhA2ko/"

# Row 8 now holds the fixture previously shown on row 6
$ws.Cells.Item(8, 6).Value = "Montpellier"
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(8, 8).Value = "Le Havre"
$ws.Cells.Item(8, 9).Value = 2
$ws.Cells.Item(8, 10).Value = 2.37
$ws.Cells.Item(8, 11).Value = "05/07/2023 22:22"
$ws.Cells.Item(8, 12).Value = 1.73
$ws.Cells.Item(8, 13).Value = "13/08/2023 14:57"
$ws.Cells.Item(8, 14).Value = 3.2
$ws.Cells.Item(8, 15).Value = "05/07/2023 22:22"
$ws.Cells.Item(8, 16).Value = 3.82
$ws.Cells.Item(8, 17).Value = "13/08/2023 14:57"
$ws.Cells.Item(8, 18).Value = 3.12
$ws.Cells.Item(8, 19).Value = "05/07/2023 22:22"
$ws.Cells.Item(8, 20).Value = 5.41
$ws.Cells.Item(8, 21).Value = "13/08/2023 14:57"
$ws.Cells.Item(8, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/montpellier-le-havre/x8rvY5sh/"

# Row 42 now holds the fixture previously shown on row 44
$ws.Cells.Item(42, 6).Value = "Clermont"
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = "Nantes"
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = 2.18
$ws.Cells.Item(42, 11).Value = "28/08/2023 16:01"
$ws.Cells.Item(42, 12).Value = 1.88
$ws.Cells.Item(42, 13).Value = "17/09/2023 14:58"
$ws.Cells.Item(42, 14).Value = 3.38
$ws.Cells.Item(42, 15).Value = "28/08/2023 16:01"
$ws.Cells.Item(42, 16).Value = 3.8
$ws.Cells.Item(42, 17).Value = "17/09/2023 14:58"
$ws.Cells.Item(42, 18).Value = 3.59
$ws.Cells.Item(42, 19).Value = "28/08/2023 16:01"
$ws.Cells.Item(42, 20).Value = 4.34
$ws.Cells.Item(42, 21).Value = "17/09/2023 14:58"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/clermont-nantes/Ox0rt4Ya/"

# Row 43 now holds the fixture previously shown on row 42
$ws.Cells.Item(43, 6).Value = "Reims"
$ws.Cells.Item(43, 7).Value = 1
$ws.Cells.Item(43, 8).Value = "Brest"
$ws.Cells.Item(43, 9).Value = 2
$ws.Cells.Item(43, 10).Value = 1.71
$ws.Cells.Item(43, 11).Value = "28/08/2023 16:01"
$ws.Cells.Item(43, 12).Value = 2.03
$ws.Cells.Item(43, 13).Value = "17/09/2023 14:50"
$ws.Cells.Item(43, 14).Value = 3.98
$ws.Cells.Item(43, 15).Value = "28/08/2023 16:01"
$ws.Cells.Item(43, 16).Value = 3.6
$ws.Cells.Item(43, 17).Value = "17/09/2023 14:53"
$ws.Cells.Item(43, 18).Value = 5.01
$ws.Cells.Item(43, 19).Value = "28/08/2023 16:01"
$ws.Cells.Item(43, 20).Value = 3.96
$ws.Cells.Item(43, 21).Value = "17/09/2023 14:53"
$ws.Cells.Item(43, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/reims-brest/pn1vspJg/"

# Row 44 now holds the fixture previously shown on row 43
$ws.Cells.Item(44, 6).Value = "Strasbourg"
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = "Montpellier"
$ws.Cells.Item(44, 9).Value = 2
$ws.Cells.Item(44, 10).Value = 2.02
$ws.Cells.Item(44, 11).Value = "28/08/2023 16:01"
$ws.Cells.Item(44, 12).Value = 3.1
$ws.Cells.Item(44, 13).Value = "17/09/2023 14:58"
$ws.Cells.Item(44, 14).Value = 3.61
$ws.Cells.Item(44, 15).Value = "28/08/2023 16:01"
$ws.Cells.Item(44, 16).Value = 3.37
$ws.Cells.Item(44, 17).Value = "17/09/2023 14:31"
$ws.Cells.Item(44, 18).Value = 3.8
$ws.Cells.Item(44, 19).Value = "28/08/2023 16:01"
$ws.Cells.Item(44, 20).Value = 2.44
$ws.Cells.Item(44, 21).Value = "17/09/2023 14:57"
$ws.Cells.Item(44, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/strasbourg-montpellier/fJq2dPIt/"

# Row 69 now holds the fixture previously shown on row 70
$ws.Cells.Item(69, 6).Value = "Lyon"
$ws.Cells.Item(69, 7).Value = 3
$ws.Cells.Item(69, 8).Value = "Lorient"
$ws.Cells.Item(69, 9).Value = 3
$ws.Cells.Item(69, 10).Value = 1.6
$ws.Cells.Item(69, 11).Value = "24/09/2023 10:02"
$ws.Cells.Item(69, 12).Value = 1.79
$ws.Cells.Item(69, 13).Value = "08/10/2023 14:55"
$ws.Cells.Item(69, 14).Value = 4.37
$ws.Cells.Item(69, 15).Value = "24/09/2023 10:02"
$ws.Cells.Item(69, 16).Value = 3.93
$ws.Cells.Item(69, 17).Value = "08/10/2023 14:58"
$ws.Cells.Item(69, 18).Value = 5.45
$ws.Cells.Item(69, 19).Value = "24/09/2023 10:02"
$ws.Cells.Item(69, 20).Value = 4.66
$ws.Cells.Item(69, 21).Value = "08/10/2023 14:57"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/lyon-lorient/Qm3i5q54/"

# Row 70 now holds the fixture previously shown on row 69
$ws.Cells.Item(70, 6).Value = "Brest"
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = "Toulouse"
$ws.Cells.Item(70, 9).Value = 1
$ws.Cells.Item(70, 10).Value = 2.36
$ws.Cells.Item(70, 11).Value = "24/09/2023 10:02"
$ws.Cells.Item(70, 12).Value = 1.73
$ws.Cells.Item(70, 13).Value = "08/10/2023 14:58"
$ws.Cells.Item(70, 14).Value = 3.47
$ws.Cells.Item(70, 15).Value = "24/09/2023 10:02"
$ws.Cells.Item(70, 16).Value = 3.9
$ws.Cells.Item(70, 17).Value = "08/10/2023 14:58"
$ws.Cells.Item(70, 18).Value = 3.11
$ws.Cells.Item(70, 19).Value = "24/09/2023 10:02"
$ws.Cells.Item(70, 20).Value = 5.25
$ws.Cells.Item(70, 21).Value = "08/10/2023 14:58"
$ws.Cells.Item(70, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/brest-toulouse/QBszcLCp/"

# Row 77 now holds the fixture previously shown on row 78
$ws.Cells.Item(77, 6).Value = "Toulouse"
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = "Reims"
$ws.Cells.Item(77, 9).Value = 1
$ws.Cells.Item(77, 10).Value = 2.67
$ws.Cells.Item(77, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(77, 12).Value = 2.73
$ws.Cells.Item(77, 13).Value = "22/10/2023 14:59"
$ws.Cells.Item(77, 14).Value = 3.5
$ws.Cells.Item(77, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(77, 16).Value = 3.4
$ws.Cells.Item(77, 17).Value = "22/10/2023 14:52"
$ws.Cells.Item(77, 18).Value = 2.53
$ws.Cells.Item(77, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(77, 20).Value = 2.76
$ws.Cells.Item(77, 21).Value = "22/10/2023 14:59"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/toulouse-reims/67KMKO4i/"

# Row 78 now holds the fixture previously shown on row 79
$ws.Cells.Item(78, 6).Value = "Nantes"
$ws.Cells.Item(78, 7).Value = 2
$ws.Cells.Item(78, 8).Value = "Montpellier"
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 2.8
$ws.Cells.Item(78, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(78, 12).Value = 3.16
$ws.Cells.Item(78, 13).Value = "22/10/2023 14:57"
$ws.Cells.Item(78, 14).Value = 3.34
$ws.Cells.Item(78, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(78, 16).Value = 3.65
$ws.Cells.Item(78, 17).Value = "22/10/2023 14:55"
$ws.Cells.Item(78, 18).Value = 2.65
$ws.Cells.Item(78, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(78, 20).Value = 2.32
$ws.Cells.Item(78, 21).Value = "22/10/2023 14:57"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/nantes-montpellier/tnvBbYmk/"

# Row 79 now holds the fixture previously shown on row 77
$ws.Cells.Item(79, 6).Value = "Lille"
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = "Brest"
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 1.52
$ws.Cells.Item(79, 11).Value = "01/10/2023 23:01"
$ws.Cells.Item(79, 12).Value = 1.91
$ws.Cells.Item(79, 13).Value = "22/10/2023 14:59"
$ws.Cells.Item(79, 14).Value = 4.48
$ws.Cells.Item(79, 15).Value = "01/10/2023 23:01"
$ws.Cells.Item(79, 16).Value = 3.55
$ws.Cells.Item(79, 17).Value = "22/10/2023 14:59"
$ws.Cells.Item(79, 18).Value = 6.48
$ws.Cells.Item(79, 19).Value = "01/10/2023 23:01"
$ws.Cells.Item(79, 20).Value = 4.61
$ws.Cells.Item(79, 21).Value = "22/10/2023 14:59"
$ws.Cells.Item(79, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/lille-brest/hOMYHMlA/"

# --- Append two newly scraped matches as rows 91 and 92 ---

# -- Row 91 --
$ws.Range("A90").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$ws.Range("E90").Copy()
$ws.Range("E91").PasteSpecial(-4122)
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = "france"
$ws.Cells.Item(91, 3).Value = "ligue-1"
$ws.Cells.Item(91, 4).Value = "2023-2024"
$ws.Cells.Item(91, 5).Value = 45234.70833333334
$ws.Cells.Item(91, 6).Value = "Lorient"
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = "Lens"
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 3.88
$ws.Cells.Item(91, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(91, 12).Value = 4.52
$ws.Cells.Item(91, 13).Value = "04/11/2023 16:59"
$ws.Cells.Item(91, 14).Value = 3.64
$ws.Cells.Item(91, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(91, 16).Value = 3.62
$ws.Cells.Item(91, 17).Value = "04/11/2023 16:59"
$ws.Cells.Item(91, 18).Value = 1.99
$ws.Cells.Item(91, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(91, 20).Value = 1.88
$ws.Cells.Item(91, 21).Value = "04/11/2023 16:54"
$ws.Cells.Item(91, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/lorient-lens/Us1sGLzK/"

# -- Row 92 --
$ws.Range("A90").Copy()
$ws.Range("A92").PasteSpecial(-4122)
$ws.Range("E90").Copy()
$ws.Range("E92").PasteSpecial(-4122)
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "france"
$ws.Cells.Item(92, 3).Value = "ligue-1"
$ws.Cells.Item(92, 4).Value = "2023-2024"
$ws.Cells.Item(92, 5).Value = 45234.875
$ws.Cells.Item(92, 6).Value = "Marseille"
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = "Lille"
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 1.94
$ws.Cells.Item(92, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(92, 12).Value = 1.93
$ws.Cells.Item(92, 13).Value = "04/11/2023 20:58"
$ws.Cells.Item(92, 14).Value = 3.67
$ws.Cells.Item(92, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(92, 16).Value = 3.7
$ws.Cells.Item(92, 17).Value = "04/11/2023 20:58"
$ws.Cells.Item(92, 18).Value = 3.71
$ws.Cells.Item(92, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(92, 20).Value = 4.17
$ws.Cells.Item(92, 21).Value = "04/11/2023 20:59"
$ws.Cells.Item(92, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/marseille-lille/ji2wH1LD/"

$excel.CutCopyMode = 0

